$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "97.085.83"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.10%  "
Set-TextValue $ws.Cells.Item(3, 4) "3.671.27"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.79%  "
Set-TextValue $ws.Cells.Item(4, 4) "2.67"
Set-TextValue $ws.Cells.Item(4, 5) "  +40.28%  "
Set-TextValue $ws.Cells.Item(5, 5) "  +0.07%  "
Set-TextValue $ws.Cells.Item(6, 4) "228.74"
Set-TextValue $ws.Cells.Item(6, 5) "  -4.28%  "
Set-TextValue $ws.Cells.Item(7, 4) "649.58"
Set-TextValue $ws.Cells.Item(7, 5) "  -1.96%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.435"
Set-TextValue $ws.Cells.Item(8, 5) "  +1.59%  "
Set-TextValue $ws.Cells.Item(9, 4) "1.21"
Set-TextValue $ws.Cells.Item(9, 5) "  +12.68%  "
Set-TextValue $ws.Cells.Item(10, 5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(11, 4) "3.670.51"
Set-TextValue $ws.Cells.Item(11, 5) "  -1.76%  "
Set-TextValue $ws.Cells.Item(12, 4) "49.26"
Set-TextValue $ws.Cells.Item(12, 5) "  +8.76%  "
Set-TextValue $ws.Cells.Item(13, 5) "  +1.76%  "
Set-TextValue $ws.Cells.Item(14, 4) "0.0000298"
Set-TextValue $ws.Cells.Item(14, 5) "  -8.13%  "
Set-TextValue $ws.Cells.Item(15, 4) "6.72"
Set-TextValue $ws.Cells.Item(15, 5) "  -3.22%  "
Set-TextValue $ws.Cells.Item(16, 4) "4.356.60"
Set-TextValue $ws.Cells.Item(16, 5) "  -1.76%  "
Set-TextValue $ws.Cells.Item(17, 4) "96.759.12"
Set-TextValue $ws.Cells.Item(17, 5) "  -0.37%  "
Set-TextValue $ws.Cells.Item(18, 4) "21.45"
Set-TextValue $ws.Cells.Item(18, 5) "  +13.45%  "
Set-TextValue $ws.Cells.Item(19, 4) "8.92"
Set-TextValue $ws.Cells.Item(19, 5) "  -2.35%  "
Set-TextValue $ws.Cells.Item(20, 4) "14.22"
Set-TextValue $ws.Cells.Item(20, 5) "  +7.86%  "
Set-TextValue $ws.Cells.Item(21, 4) "3.667.80"
Set-TextValue $ws.Cells.Item(21, 5) "  -1.64%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.557"
Set-TextValue $ws.Cells.Item(22, 5) "  +9.48%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.276"
Set-TextValue $ws.Cells.Item(23, 5) "  +41.99%  "
Set-TextValue $ws.Cells.Item(24, 4) "529.55"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.14%  "
Set-TextValue $ws.Cells.Item(25, 4) "3.30"
Set-TextValue $ws.Cells.Item(25, 5) "  -5.75%  "
Set-TextValue $ws.Cells.Item(26, 4) "122.69"
Set-TextValue $ws.Cells.Item(26, 5) "  +12.20%  "
Set-TextValue $ws.Cells.Item(27, 4) "0.0000207"
Set-TextValue $ws.Cells.Item(27, 5) "  -8.03%  "
Set-TextValue $ws.Cells.Item(28, 4) "6.87"
Set-TextValue $ws.Cells.Item(28, 5) "  -0.97%  "
Set-TextValue $ws.Cells.Item(29, 4) "3.847.17"
Set-TextValue $ws.Cells.Item(29, 5) "  -1.73%  "
Set-TextValue $ws.Cells.Item(30, 4) "13.13"
Set-TextValue $ws.Cells.Item(30, 5) "  -5.14%  "
Set-TextValue $ws.Cells.Item(31, 5) "  +0.76%  "
Set-TextValue $ws.Cells.Item(32, 4) "3.03"
Set-TextValue $ws.Cells.Item(32, 5) "  -0.94%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.999"
Set-TextValue $ws.Cells.Item(33, 5) "  -0.18%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.185"
Set-TextValue $ws.Cells.Item(34, 5) "  -4.03%  "
Set-TextValue $ws.Cells.Item(35, 4) "33.22"
Set-TextValue $ws.Cells.Item(35, 5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.623"
Set-TextValue $ws.Cells.Item(36, 5) "  +4.03%  "
Set-TextValue $ws.Cells.Item(37, 5) "  +0.54%  "
Set-TextValue $ws.Cells.Item(38, 4) "1.79"
Set-TextValue $ws.Cells.Item(38, 5) "  -3.96%  "
Set-TextValue $ws.Cells.Item(39, 4) "607.92"
Set-TextValue $ws.Cells.Item(39, 5) "  -6.81%  "
Set-TextValue $ws.Cells.Item(40, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(41, 4) "8.56"
Set-TextValue $ws.Cells.Item(41, 5) "  -3.94%  "
Set-TextValue $ws.Cells.Item(42, 4) "7.11"
Set-TextValue $ws.Cells.Item(42, 5) "  +3.07%  "
Set-TextValue $ws.Cells.Item(43, 4) "42.36"
Set-TextValue $ws.Cells.Item(43, 5) "  +1.54%  "
Set-TextValue $ws.Cells.Item(44, 5) "  +9.74%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.493"
Set-TextValue $ws.Cells.Item(45, 5) "  +2.88%  "
Set-TextValue $ws.Cells.Item(46, 5) "  -5.44%  "
Set-TextValue $ws.Cells.Item(47, 5) "  -2.04%  "
Set-TextValue $ws.Cells.Item(48, 4) "1.97"
Set-TextValue $ws.Cells.Item(48, 5) "  -4.14%  "
Set-TextValue $ws.Cells.Item(49, 4) "236.17"
Set-TextValue $ws.Cells.Item(49, 5) "  +13.62%  "
Set-TextValue $ws.Cells.Item(50, 5) "  -3.86%  "
Set-TextValue $ws.Cells.Item(51, 4) "8.93"
Set-TextValue $ws.Cells.Item(51, 5) "  +1.54%  "
